$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.102.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6162"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07331"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2891"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07659"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.820.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.945"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6605"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "81.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008952"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.841"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.090.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.055.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.125"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1405"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.425"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.480"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05578"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.088"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.103"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.818"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7337"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.622"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.836"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.204.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.374"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8915"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.959.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000120"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5084"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.076"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05788"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.68%  "
